# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# This updates the recalibrated ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN columns (D:H)
# for the affected rows in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> D,E,F,G,H values
$updates = @{
    3  = @(118145.4429732955, -0.03610817978583215, 0.1919615532859447, -1.211843204895535, 11.82180579330556)
    4  = @(118884.566777737,  -0.04134920333870538, 0.2341443669460976, -1.554505504497392, 13.29899062754492)
    5  = @(120012.6945115529, -0.03198243896603334, 0.2111750862737707, -0.5194278457725386, 6.391204359380751)
    6  = @(120339.5382767707, -0.04724186125118832, 0.2612644695633924, -1.270849681707819,  10.89516532113534)
    7  = @(121825.6279993976, -0.04355950020084948, 0.2470381920871975, -0.5372328021474828, 4.607595684585159)
    8  = @(122013.8068534817, -0.05885547107895549, 0.2247380250754349, -0.968356254154783,  7.193524516484642)
    9  = @(123578.8200288022, -0.08842387603249097, 0.3543224698284816, -1.943705335092959,  12.41810975826059)
    10 = @(124980.4375542639, -0.1245256475804394,  0.4433245553660382, -1.931558653852689,  9.848049547071152)
    11 = @(127045.4898404461, -0.2012459396652004,  0.7802362637687535, -2.622655670293256,  12.79669364474417)
    17 = @(117475.4294382372, -0.04786276384093117, 0.1235721125547685, -1.27235776433512,   13.1051602762496)
    18 = @(117483.2433203897, -0.04121138019182476, 0.132879978150446,  -0.875721473036936,  10.1014731699809)
    19 = @(117501.1614142006, -0.04052167028609403, 0.1419128143011636, -0.5826305039107675, 6.426602499682184)
    20 = @(118203.7846888634, -0.0266863524811118,  0.1808508380883164, -1.493562302033354,  13.06899994247547)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
